$wb = $excel.ActiveWorkbook

# AR sheet
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.03005027436747393
$ws.Range("B3").Value = 0.7326361358633129
$ws.Range("B4").Value = 0.1409879752221953

# SETAR sheet
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B2").Value = -0.1975709897101832
$ws.Range("B3").Value = 0.4713572626792478
$ws.Range("B4").Value = 0.08614862325455526
$ws.Range("B5").Value = 0.2213492672759159
$ws.Range("B6").Value = 0.5201949097290097
$ws.Range("B7").Value = 0.1108289563811776

# GARCH sheet
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.001720135121303271
$ws.Range("B3").Value = 0.1241085095164853
$ws.Range("B4").Value = 0.1170763193941317
$ws.Range("B5").Value = 0.09727026973382057

# TARCH sheet
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.005194563981180916
$ws.Range("B3").Value = 0.1270243165423624
$ws.Range("B4").Value = 0.1624317820460614
$ws.Range("B5").Value = -0.09565228952854361
$ws.Range("B6").Value = 0.08052319728702199

# AR-TARCH sheet
$ws = $wb.Worksheets.Item("AR-TARCH")
$ws.Range("B2").Value = 0.02994676599299398
$ws.Range("B3").Value = 0.742207539510004
$ws.Range("B4").Value = 0.1226048382404183
$ws.Range("B5").Value = 0.1601747738778287
$ws.Range("B6").Value = -0.08788668274862838
$ws.Range("B7").Value = 267.2235201090382 * 0.000000000000001
